# Updates the cryptos list worksheet with refreshed price / volume data
# (GitHub Actions data refresh), fixes the ordering of the "Stacks" /
# "Cosmos" rows (45/46), and replaces the "EnergySwap" row (51) with
# "Arweave" data, matching the upstream commit.

# Helper: force a numeric-looking string to be stored as TEXT (not a
# number) in the cell, mirroring how the source data is authored, then
# strip the temporary "Text" number-format so the cell is left with no
# extra style applied (same as its original, un-styled state).
function Set-TextValue {
    param($Sheet, $Address, $Text)
    $rng = $Sheet.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.558.04'
$ws.Range("E2").Value = '  -0.37%  '

$ws.Range("D3").Value = '3.725.92'
$ws.Range("E3").Value = '  -1.92%  '

$ws.Range("E4").Value = '  -0.15%  '

Set-TextValue $ws "D5" '591.26'
$ws.Range("E5").Value = '  -1.14%  '

Set-TextValue $ws "D6" '164.93'
$ws.Range("E6").Value = '  -2.32%  '

$ws.Range("D7").Value = '3.726.58'

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  -2.09%  '

$ws.Range("E10").Value = '  -4.00%  '

Set-TextValue $ws "D11" '6.41'
$ws.Range("E11").Value = '  -1.02%  '

$ws.Range("E13").Value = '  -5.77%  '

Set-TextValue $ws "D14" '35.88'
$ws.Range("E14").Value = '  -2.42%  '

$ws.Range("D15").Value = '4.348.94'
$ws.Range("E15").Value = '  -1.91%  '

$ws.Range("D16").Value = '3.722.66'
$ws.Range("E16").Value = '  -1.35%  '

$ws.Range("D17").Value = '67.522.91'
$ws.Range("E17").Value = '  -0.54%  '

$ws.Range("E18").Value = '  +0.47%  '

Set-TextValue $ws "D19" '7.03'
$ws.Range("E19").Value = '  -5.15%  '

$ws.Range("E20").Value = '  -0.35%  '

$ws.Range("E21").Value = '  -2.11%  '

Set-TextValue $ws "D22" '464.24'
$ws.Range("E22").Value = '  -1.05%  '

$ws.Range("E23").Value = '  -3.69%  '

Set-TextValue $ws "D24" '82.49'
$ws.Range("E24").Value = '  -0.86%  '

Set-TextValue $ws "D25" '0.0000135'
$ws.Range("E25").Value = '  -10.00%  '

Set-TextValue $ws "D26" '2.17'
$ws.Range("E26").Value = '  -4.10%  '

$ws.Range("E27").Value = '  -1.79%  '

Set-TextValue $ws "D28" '10.13'
$ws.Range("E28").Value = '  -1.24%  '

$ws.Range("D30").Value = '3.870.16'
$ws.Range("E30").Value = '  -1.93%  '

$ws.Range("E31").Value = '  -6.00%  '

Set-TextValue $ws "D32" '7.33'
$ws.Range("E32").Value = '  -5.05%  '

$ws.Range("E33").Value = '  -3.38%  '

$ws.Range("E34").Value = '  -3.83%  '

Set-TextValue $ws "D35" '8.97'
$ws.Range("E35").Value = '  -3.77%  '

$ws.Range("D36").Value = '3.676.59'
$ws.Range("E36").Value = '  -2.30%  '

$ws.Range("E37").Value = '  -5.43%  '

$ws.Range("E38").Value = '  -10.49%  '

Set-TextValue $ws "D39" '0.136'
$ws.Range("E39").Value = '  -2.39%  '

Set-TextValue $ws "D40" '0.989'
$ws.Range("E40").Value = '  -2.58%  '

Set-TextValue $ws "D41" '5.72'
$ws.Range("E41").Value = '  -3.74%  '

$ws.Range("E42").Value = '  -0.11%  '

$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("E44").Value = '  -3.57%  '

$ws.Range("B45").Value = 'Cosmos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws "D45" '8.49'
$ws.Range("E45").Value = '  -3.98%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws "D46" '1.91'
$ws.Range("E46").Value = '  -3.20%  '

Set-TextValue $ws "D47" '45.33'
$ws.Range("E47").Value = '  -2.62%  '

Set-TextValue $ws "D48" '392.53'
$ws.Range("E48").Value = '  -3.70%  '

Set-TextValue $ws "D49" '143.36'
$ws.Range("E49").Value = '  +1.42%  '

$ws.Range("E50").Value = '  -3.31%  '

$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws "D51" '38.32'
$ws.Range("E51").Value = '  +1.17%  '
